$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in H1, copying the formatting (bold, borders,
# centered alignment) from the existing header cell G1 ("sum")
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Add the corresponding data value in H2 (default/unstyled, like the
# other numeric cells in row 2)
$ws.Range("H2").Value = 1
